$p = $ppt.ActivePresentation

# --- Slide 2: merge "Ed " + "Salinas" runs into a single run "Ed Salinas" ---
$s2 = $p.Slides.Item(2)
$subtitle = $s2.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Paragraphs(1).Text = "Ed Salinas"

# --- Slide 2: TextBox 2 -> update "tayo" to "sa", "makipagdate" to "lablyf " ---
$textBox2 = $s2.Shapes.Item(3)
$tr = $textBox2.TextFrame.TextRange
$tr.Runs(2).Text = "sa"
$tr.Runs(4).Text = "lablyf "

# --- Slide 22: reposition/resize the picture ---
$s22 = $p.Slides.Item(22)
$pic = $s22.Shapes.Item(2)
$pic.Left = 2894012 / 12700
$pic.Top = 1071563 / 12700
$pic.Width = 5410200 / 12700
$pic.Height = 5607123 / 12700
